$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("Tidur") replaces row 5 ("Makan"): read row 6's values first...
$a6 = $ws.Range("A6").Value2
$b6 = $ws.Range("B6").Value2
$c6 = $ws.Range("C6").Value2
$d6 = $ws.Range("D6").Value2
$e6 = $ws.Range("E6").Value2

# ...then write them into row 5, overwriting the old "Makan" entry.
$ws.Range("A5").Value = $a6
$ws.Range("B5").Value = $b6
$ws.Range("C5").Value = $c6
$ws.Range("D5").Value = $d6
$ws.Range("E5").Value = $e6

# Remove the now-duplicated row 6 (rows below stay put; no shift).
$null = $ws.Range("A6:E6").ClearContents()

# Match the resulting selection recorded in the sheet view.
$null = $ws.Range("A5:E5").Select()
